$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the "time_taken" timestamps in column F of the data sheet
$dataSheet.Range("F2").Value = "2021-10-05 14:22:52.889667"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:52.889674"
$dataSheet.Range("F4").Value = "2021-10-05 14:22:52.889678"
$dataSheet.Range("F5").Value = "2021-10-05 14:22:52.889681"
$dataSheet.Range("F6").Value = "2021-10-05 14:22:52.889683"
$dataSheet.Range("F7").Value = "2021-10-05 14:22:52.889686"
$dataSheet.Range("F8").Value = "2021-10-05 14:22:52.889689"
$dataSheet.Range("F9").Value = "2021-10-05 14:22:52.889691"
$dataSheet.Range("F10").Value = "2021-10-05 14:22:52.889694"
$dataSheet.Range("F11").Value = "2021-10-05 14:22:52.889697"
$dataSheet.Range("F12").Value = "2021-10-05 14:22:52.889699"
$dataSheet.Range("F13").Value = "2021-10-05 14:22:52.889702"
$dataSheet.Range("F14").Value = "2021-10-05 14:22:52.889704"
$dataSheet.Range("F15").Value = "2021-10-05 14:22:52.889707"
$dataSheet.Range("F16").Value = "2021-10-05 14:22:52.889710"
$dataSheet.Range("F17").Value = "2021-10-05 14:22:52.889712"
$dataSheet.Range("F18").Value = "2021-10-05 14:22:52.889715"
$dataSheet.Range("F19").Value = "2021-10-05 14:22:52.889718"
$dataSheet.Range("F20").Value = "2021-10-05 14:22:52.889720"
$dataSheet.Range("F21").Value = "2021-10-05 14:22:52.889723"
$dataSheet.Range("F22").Value = "2021-10-05 14:22:52.889725"
$dataSheet.Range("F23").Value = "2021-10-05 14:22:52.889728"
$dataSheet.Range("F24").Value = "2021-10-05 14:22:52.889730"
$dataSheet.Range("F25").Value = "2021-10-05 14:22:52.889732"
$dataSheet.Range("F26").Value = "2021-10-05 14:22:52.889735"
$dataSheet.Range("F27").Value = "2021-10-05 14:22:52.889738"
$dataSheet.Range("F28").Value = "2021-10-05 14:22:52.889740"
$dataSheet.Range("F29").Value = "2021-10-05 14:22:52.889743"
$dataSheet.Range("F30").Value = "2021-10-05 14:22:52.889745"
$dataSheet.Range("F31").Value = "2021-10-05 14:22:52.889748"
$dataSheet.Range("F32").Value = "2021-10-05 14:22:52.889750"
$dataSheet.Range("F33").Value = "2021-10-05 14:22:52.889753"
$dataSheet.Range("F34").Value = "2021-10-05 14:22:52.889755"
$dataSheet.Range("F35").Value = "2021-10-05 14:22:52.889758"
$dataSheet.Range("F36").Value = "2021-10-05 14:22:52.889761"
$dataSheet.Range("F37").Value = "2021-10-05 14:22:52.889763"
$dataSheet.Range("F38").Value = "2021-10-05 14:22:52.889765"
$dataSheet.Range("F39").Value = "2021-10-05 14:22:52.889768"
$dataSheet.Range("F40").Value = "2021-10-05 14:22:52.889771"
$dataSheet.Range("F41").Value = "2021-10-05 14:22:52.889773"
$dataSheet.Range("F42").Value = "2021-10-05 14:22:52.889776"
$dataSheet.Range("F43").Value = "2021-10-05 14:22:52.889779"
$dataSheet.Range("F44").Value = "2021-10-05 14:22:52.889781"
$dataSheet.Range("F45").Value = "2021-10-05 14:22:52.889784"
$dataSheet.Range("F46").Value = "2021-10-05 14:22:52.889786"
$dataSheet.Range("F47").Value = "2021-10-05 14:22:52.889789"
$dataSheet.Range("F48").Value = "2021-10-05 14:22:52.889791"
$dataSheet.Range("F49").Value = "2021-10-05 14:22:52.889794"
$dataSheet.Range("F50").Value = "2021-10-05 14:22:52.889796"
$dataSheet.Range("F51").Value = "2021-10-05 14:22:52.889799"
$dataSheet.Range("F52").Value = "2021-10-05 14:22:52.889801"
$dataSheet.Range("F53").Value = "2021-10-05 14:22:52.889804"
$dataSheet.Range("F54").Value = "2021-10-05 14:22:52.889806"
$dataSheet.Range("F55").Value = "2021-10-05 14:22:52.889809"
$dataSheet.Range("F56").Value = "2021-10-05 14:22:52.889812"
$dataSheet.Range("F57").Value = "2021-10-05 14:22:52.889814"
$dataSheet.Range("F58").Value = "2021-10-05 14:22:52.889817"
$dataSheet.Range("F59").Value = "2021-10-05 14:22:52.889819"
$dataSheet.Range("F60").Value = "2021-10-05 14:22:52.889821"
$dataSheet.Range("F61").Value = "2021-10-05 14:22:52.889824"
$dataSheet.Range("F62").Value = "2021-10-05 14:22:52.889826"
$dataSheet.Range("F63").Value = "2021-10-05 14:22:52.889829"
$dataSheet.Range("F64").Value = "2021-10-05 14:22:52.889831"
$dataSheet.Range("F65").Value = "2021-10-05 14:22:52.889834"
$dataSheet.Range("F66").Value = "2021-10-05 14:22:52.889837"
$dataSheet.Range("F67").Value = "2021-10-05 14:22:52.889840"
$dataSheet.Range("F68").Value = "2021-10-05 14:22:52.889843"
$dataSheet.Range("F69").Value = "2021-10-05 14:22:52.889845"
$dataSheet.Range("F70").Value = "2021-10-05 14:22:52.889847"
$dataSheet.Range("F71").Value = "2021-10-05 14:22:52.889850"
$dataSheet.Range("F72").Value = "2021-10-05 14:22:52.889852"
$dataSheet.Range("F73").Value = "2021-10-05 14:22:52.889855"
$dataSheet.Range("F74").Value = "2021-10-05 14:22:52.889858"
$dataSheet.Range("F75").Value = "2021-10-05 14:22:52.889860"
$dataSheet.Range("F76").Value = "2021-10-05 14:22:52.889863"
$dataSheet.Range("F77").Value = "2021-10-05 14:22:52.889865"
$dataSheet.Range("F78").Value = "2021-10-05 14:22:52.889870"
$dataSheet.Range("F79").Value = "2021-10-05 14:22:52.889873"
$dataSheet.Range("F80").Value = "2021-10-05 14:22:52.889875"
$dataSheet.Range("F81").Value = "2021-10-05 14:22:52.889877"
$dataSheet.Range("F82").Value = "2021-10-05 14:22:52.889880"
$dataSheet.Range("F83").Value = "2021-10-05 14:22:52.889883"
$dataSheet.Range("F84").Value = "2021-10-05 14:22:52.889885"
$dataSheet.Range("F85").Value = "2021-10-05 14:22:52.889887"
$dataSheet.Range("F86").Value = "2021-10-05 14:22:52.889890"
$dataSheet.Range("F87").Value = "2021-10-05 14:22:52.889892"
$dataSheet.Range("F88").Value = "2021-10-05 14:22:52.889895"
$dataSheet.Range("F89").Value = "2021-10-05 14:22:52.889897"
$dataSheet.Range("F90").Value = "2021-10-05 14:22:52.889900"
$dataSheet.Range("F91").Value = "2021-10-05 14:22:52.889902"
$dataSheet.Range("F92").Value = "2021-10-05 14:22:52.889905"
$dataSheet.Range("F93").Value = "2021-10-05 14:22:52.889907"
$dataSheet.Range("F94").Value = "2021-10-05 14:22:52.889911"
$dataSheet.Range("F95").Value = "2021-10-05 14:22:52.889914"
$dataSheet.Range("F96").Value = "2021-10-05 14:22:52.889916"
$dataSheet.Range("F97").Value = "2021-10-05 14:22:52.889919"
$dataSheet.Range("F98").Value = "2021-10-05 14:22:52.889922"
$dataSheet.Range("F99").Value = "2021-10-05 14:22:52.889924"
$dataSheet.Range("F100").Value = "2021-10-05 14:22:52.889927"
$dataSheet.Range("F101").Value = "2021-10-05 14:22:52.889929"
$dataSheet.Range("F102").Value = "2021-10-05 14:22:52.889932"
$dataSheet.Range("F103").Value = "2021-10-05 14:22:52.889934"
$dataSheet.Range("F104").Value = "2021-10-05 14:22:52.889937"
$dataSheet.Range("F105").Value = "2021-10-05 14:22:52.889939"
$dataSheet.Range("F106").Value = "2021-10-05 14:22:52.889942"
$dataSheet.Range("F107").Value = "2021-10-05 14:22:52.889944"
$dataSheet.Range("F108").Value = "2021-10-05 14:22:52.889947"
$dataSheet.Range("F109").Value = "2021-10-05 14:22:52.889949"
$dataSheet.Range("F110").Value = "2021-10-05 14:22:52.889954"
$dataSheet.Range("F111").Value = "2021-10-05 14:22:52.889957"
$dataSheet.Range("F112").Value = "2021-10-05 14:22:52.889959"
$dataSheet.Range("F113").Value = "2021-10-05 14:22:52.889962"
$dataSheet.Range("F114").Value = "2021-10-05 14:22:52.889965"
$dataSheet.Range("F115").Value = "2021-10-05 14:22:52.889967"
$dataSheet.Range("F116").Value = "2021-10-05 14:22:52.889970"
$dataSheet.Range("F117").Value = "2021-10-05 14:22:52.889973"
$dataSheet.Range("F118").Value = "2021-10-05 14:22:52.889975"
$dataSheet.Range("F119").Value = "2021-10-05 14:22:52.889978"
$dataSheet.Range("F120").Value = "2021-10-05 14:22:52.889980"
$dataSheet.Range("F121").Value = "2021-10-05 14:22:52.889983"
$dataSheet.Range("F122").Value = "2021-10-05 14:22:52.889986"
$dataSheet.Range("F123").Value = "2021-10-05 14:22:52.889988"
$dataSheet.Range("F124").Value = "2021-10-05 14:22:52.889991"
$dataSheet.Range("F125").Value = "2021-10-05 14:22:52.889993"
$dataSheet.Range("F126").Value = "2021-10-05 14:22:52.889996"
$dataSheet.Range("F127").Value = "2021-10-05 14:22:52.889998"
$dataSheet.Range("F128").Value = "2021-10-05 14:22:52.890001"
$dataSheet.Range("F129").Value = "2021-10-05 14:22:52.890004"
$dataSheet.Range("F130").Value = "2021-10-05 14:22:52.890008"
$dataSheet.Range("F131").Value = "2021-10-05 14:22:52.890011"
$dataSheet.Range("F132").Value = "2021-10-05 14:22:52.890014"
$dataSheet.Range("F133").Value = "2021-10-05 14:22:52.890017"
$dataSheet.Range("F134").Value = "2021-10-05 14:22:52.890019"
$dataSheet.Range("F135").Value = "2021-10-05 14:22:52.890022"
$dataSheet.Range("F136").Value = "2021-10-05 14:22:52.890024"
$dataSheet.Range("F137").Value = "2021-10-05 14:22:52.890027"
$dataSheet.Range("F138").Value = "2021-10-05 14:22:52.890030"
$dataSheet.Range("F139").Value = "2021-10-05 14:22:52.890032"
$dataSheet.Range("F140").Value = "2021-10-05 14:22:52.890035"
$dataSheet.Range("F141").Value = "2021-10-05 14:22:52.890037"
$dataSheet.Range("F142").Value = "2021-10-05 14:22:52.890040"
$dataSheet.Range("F143").Value = "2021-10-05 14:22:52.890042"
$dataSheet.Range("F144").Value = "2021-10-05 14:22:52.890045"
$dataSheet.Range("F145").Value = "2021-10-05 14:22:52.890048"
$dataSheet.Range("F146").Value = "2021-10-05 14:22:52.890050"
$dataSheet.Range("F147").Value = "2021-10-05 14:22:52.890053"
$dataSheet.Range("F148").Value = "2021-10-05 14:22:52.890055"
$dataSheet.Range("F149").Value = "2021-10-05 14:22:52.890058"
$dataSheet.Range("F150").Value = "2021-10-05 14:22:52.890061"
$dataSheet.Range("F151").Value = "2021-10-05 14:22:52.890063"
$dataSheet.Range("F152").Value = "2021-10-05 14:22:52.890066"
$dataSheet.Range("F153").Value = "2021-10-05 14:22:52.890069"
$dataSheet.Range("F154").Value = "2021-10-05 14:22:52.890071"
$dataSheet.Range("F155").Value = "2021-10-05 14:22:52.890074"
$dataSheet.Range("F156").Value = "2021-10-05 14:22:52.890076"
$dataSheet.Range("F157").Value = "2021-10-05 14:22:52.890079"
$dataSheet.Range("F158").Value = "2021-10-05 14:22:52.890081"
$dataSheet.Range("F159").Value = "2021-10-05 14:22:52.890084"
$dataSheet.Range("F160").Value = "2021-10-05 14:22:52.890087"
$dataSheet.Range("F161").Value = "2021-10-05 14:22:52.890089"
$dataSheet.Range("F162").Value = "2021-10-05 14:22:52.890092"
$dataSheet.Range("F163").Value = "2021-10-05 14:22:52.890094"
$dataSheet.Range("F164").Value = "2021-10-05 14:22:52.890097"
$dataSheet.Range("F165").Value = "2021-10-05 14:22:52.890099"
$dataSheet.Range("F166").Value = "2021-10-05 14:22:52.890102"
$dataSheet.Range("F167").Value = "2021-10-05 14:22:52.890104"
$dataSheet.Range("F168").Value = "2021-10-05 14:22:52.890107"
$dataSheet.Range("F169").Value = "2021-10-05 14:22:52.890110"
$dataSheet.Range("F170").Value = "2021-10-05 14:22:52.890112"
$dataSheet.Range("F171").Value = "2021-10-05 14:22:52.890115"
$dataSheet.Range("F172").Value = "2021-10-05 14:22:52.890118"
$dataSheet.Range("F173").Value = "2021-10-05 14:22:52.890120"
$dataSheet.Range("F174").Value = "2021-10-05 14:22:52.890124"
$dataSheet.Range("F175").Value = "2021-10-05 14:22:52.890127"
$dataSheet.Range("F176").Value = "2021-10-05 14:22:52.890130"
$dataSheet.Range("F177").Value = "2021-10-05 14:22:52.890132"
$dataSheet.Range("F178").Value = "2021-10-05 14:22:52.890135"
$dataSheet.Range("F179").Value = "2021-10-05 14:22:52.890138"
$dataSheet.Range("F180").Value = "2021-10-05 14:22:52.890140"
$dataSheet.Range("F181").Value = "2021-10-05 14:22:52.890143"
$dataSheet.Range("F182").Value = "2021-10-05 14:22:52.890146"
$dataSheet.Range("F183").Value = "2021-10-05 14:22:52.890148"
$dataSheet.Range("F184").Value = "2021-10-05 14:22:52.890151"
$dataSheet.Range("F185").Value = "2021-10-05 14:22:52.890153"
$dataSheet.Range("F186").Value = "2021-10-05 14:22:52.890156"
$dataSheet.Range("F187").Value = "2021-10-05 14:22:52.890158"
$dataSheet.Range("F188").Value = "2021-10-05 14:22:52.890161"
$dataSheet.Range("F189").Value = "2021-10-05 14:22:52.890164"
$dataSheet.Range("F190").Value = "2021-10-05 14:22:52.890166"
$dataSheet.Range("F191").Value = "2021-10-05 14:22:52.890169"
$dataSheet.Range("F192").Value = "2021-10-05 14:22:52.890172"
$dataSheet.Range("F193").Value = "2021-10-05 14:22:52.890174"
$dataSheet.Range("F194").Value = "2021-10-05 14:22:52.890177"
$dataSheet.Range("F195").Value = "2021-10-05 14:22:52.890180"
$dataSheet.Range("F196").Value = "2021-10-05 14:22:52.890182"
$dataSheet.Range("F197").Value = "2021-10-05 14:22:52.890185"
$dataSheet.Range("F198").Value = "2021-10-05 14:22:52.890187"
$dataSheet.Range("F199").Value = "2021-10-05 14:22:52.890190"
$dataSheet.Range("F200").Value = "2021-10-05 14:22:52.890192"
$dataSheet.Range("F201").Value = "2021-10-05 14:22:52.890195"
$dataSheet.Range("F202").Value = "2021-10-05 14:22:52.890198"
$dataSheet.Range("F203").Value = "2021-10-05 14:22:52.890201"
$dataSheet.Range("F204").Value = "2021-10-05 14:22:52.890203"
$dataSheet.Range("F205").Value = "2021-10-05 14:22:52.890206"
$dataSheet.Range("F206").Value = "2021-10-05 14:22:52.890208"
$dataSheet.Range("F207").Value = "2021-10-05 14:22:52.890210"
$dataSheet.Range("F208").Value = "2021-10-05 14:22:52.890213"
$dataSheet.Range("F209").Value = "2021-10-05 14:22:52.890216"
$dataSheet.Range("F210").Value = "2021-10-05 14:22:52.890219"
$dataSheet.Range("F211").Value = "2021-10-05 14:22:52.890221"
$dataSheet.Range("F212").Value = "2021-10-05 14:22:52.890224"
$dataSheet.Range("F213").Value = "2021-10-05 14:22:52.890226"
$dataSheet.Range("F214").Value = "2021-10-05 14:22:52.890229"
$dataSheet.Range("F215").Value = "2021-10-05 14:22:52.890232"
$dataSheet.Range("F216").Value = "2021-10-05 14:22:52.890234"
$dataSheet.Range("F217").Value = "2021-10-05 14:22:52.890237"
$dataSheet.Range("F218").Value = "2021-10-05 14:22:52.890239"
$dataSheet.Range("F219").Value = "2021-10-05 14:22:52.890242"
$dataSheet.Range("F220").Value = "2021-10-05 14:22:52.890245"
$dataSheet.Range("F221").Value = "2021-10-05 14:22:52.890247"
$dataSheet.Range("F222").Value = "2021-10-05 14:22:52.890250"
$dataSheet.Range("F223").Value = "2021-10-05 14:22:52.890253"
$dataSheet.Range("F224").Value = "2021-10-05 14:22:52.890255"
$dataSheet.Range("F225").Value = "2021-10-05 14:22:52.890258"
$dataSheet.Range("F226").Value = "2021-10-05 14:22:52.890261"
$dataSheet.Range("F227").Value = "2021-10-05 14:22:52.890263"
$dataSheet.Range("F228").Value = "2021-10-05 14:22:52.890266"
$dataSheet.Range("F229").Value = "2021-10-05 14:22:52.890269"
$dataSheet.Range("F230").Value = "2021-10-05 14:22:52.890271"
$dataSheet.Range("F231").Value = "2021-10-05 14:22:52.890274"
$dataSheet.Range("F232").Value = "2021-10-05 14:22:52.890277"
$dataSheet.Range("F233").Value = "2021-10-05 14:22:52.890279"
$dataSheet.Range("F234").Value = "2021-10-05 14:22:52.890283"
$dataSheet.Range("F235").Value = "2021-10-05 14:22:52.890286"
$dataSheet.Range("F236").Value = "2021-10-05 14:22:52.890289"
$dataSheet.Range("F237").Value = "2021-10-05 14:22:52.890291"
$dataSheet.Range("F238").Value = "2021-10-05 14:22:52.890294"
$dataSheet.Range("F239").Value = "2021-10-05 14:22:52.890297"
$dataSheet.Range("F240").Value = "2021-10-05 14:22:52.890299"
$dataSheet.Range("F241").Value = "2021-10-05 14:22:52.890302"
$dataSheet.Range("F242").Value = "2021-10-05 14:22:52.890305"
$dataSheet.Range("F243").Value = "2021-10-05 14:22:52.890307"
$dataSheet.Range("F244").Value = "2021-10-05 14:22:52.890310"
$dataSheet.Range("F245").Value = "2021-10-05 14:22:52.890313"
$dataSheet.Range("F246").Value = "2021-10-05 14:22:52.890315"
$dataSheet.Range("F247").Value = "2021-10-05 14:22:52.890318"
$dataSheet.Range("F248").Value = "2021-10-05 14:22:52.890320"
$dataSheet.Range("F249").Value = "2021-10-05 14:22:52.890323"
$dataSheet.Range("F250").Value = "2021-10-05 14:22:52.890326"
$dataSheet.Range("F251").Value = "2021-10-05 14:22:52.890328"
$dataSheet.Range("F252").Value = "2021-10-05 14:22:52.890331"
$dataSheet.Range("F253").Value = "2021-10-05 14:22:52.890334"
$dataSheet.Range("F254").Value = "2021-10-05 14:22:52.890336"
$dataSheet.Range("F255").Value = "2021-10-05 14:22:52.890339"
$dataSheet.Range("F256").Value = "2021-10-05 14:22:52.890341"
$dataSheet.Range("F257").Value = "2021-10-05 14:22:52.890344"
$dataSheet.Range("F258").Value = "2021-10-05 14:22:52.890346"
$dataSheet.Range("F259").Value = "2021-10-05 14:22:52.890349"
$dataSheet.Range("F260").Value = "2021-10-05 14:22:52.890351"
$dataSheet.Range("F261").Value = "2021-10-05 14:22:52.890354"
$dataSheet.Range("F262").Value = "2021-10-05 14:22:52.890357"
$dataSheet.Range("F263").Value = "2021-10-05 14:22:52.890360"
$dataSheet.Range("F264").Value = "2021-10-05 14:22:52.890362"
$dataSheet.Range("F265").Value = "2021-10-05 14:22:52.890365"
$dataSheet.Range("F266").Value = "2021-10-05 14:22:52.890367"
$dataSheet.Range("F267").Value = "2021-10-05 14:22:52.890370"
$dataSheet.Range("F268").Value = "2021-10-05 14:22:52.890373"
# Add a new "metadata" worksheet after the "data" sheet
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1) -- bold, bordered, centered like the "data" sheet's header style
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row (row 2)
$a2 = $metaSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$metaSheet.Range("B2").Value = "Sudden cardiac death"
$metaSheet.Range("C2").Value = 841

$d2 = $metaSheet.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "9.193"
$d2.ClearFormats()

$metaSheet.Range("E2").Value = "2021-09-28T09:49:50.576191Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:52.886485"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/841/?format=json"

Write-Output "metadata sheet created"
